$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Fill in the new row of data (row 9) -- finishing the smells-analyze table
$ws.Range("A9").Value = "WebSecurityConfig.java"
$ws.Range("B9").Value = "Unutilized Abstraction: The tool detected the smell in this class because this class is potentially unused. (Please ignore the smell if the reported class is auto-generated and/or used to serve a specific known purpose.)"
$ws.Range("C9").Value = "No"
$ws.Range("D9").Value = "Checking the sercurity authetic heeaders in every .authetic request."

# Match the formatting of the other data rows (centered horizontally/vertically, wrap
# text) by copying the format from an already-styled row instead of touching each
# alignment property individually (keeps the shared cellXfs entry instead of minting
# new combinations).
$ws.Range("A6:D6").Copy()
$ws.Range("A9:D9").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Rows.Item(9).RowHeight = 126

# Update the view: scroll down a couple of rows and leave the last-typed cell selected
$excel.ActiveWindow.ScrollRow = 5
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("D9").Select()

# Window got taller in the saved view state
$excel.ActiveWindow.Height = 1102
